$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 274, shifting rows 274:334 down to 275:335
$ws.Rows.Item(274).Insert()

# Populate the new row 274 with data
$ws.Cells.Item(274, 1).Value = 10
$ws.Cells.Item(274, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(274, 3).Value = "La Araucanía"
$ws.Cells.Item(274, 4).Value = 44508
$ws.Cells.Item(274, 5).Value = 9
$ws.Cells.Item(274, 6).Value = 100112028
$ws.Cells.Item(274, 7).Value = "Sandia"
$ws.Cells.Item(274, 8).Value = "Sin especificar"
$ws.Cells.Item(274, 9).Value = "Primera"
$ws.Cells.Item(274, 10).Value = 400
$ws.Cells.Item(274, 11).Value = 800
$ws.Cells.Item(274, 12).Value = 900
$ws.Cells.Item(274, 13).Value = 850
$ws.Cells.Item(274, 14).Value = "`$/kilo (volumen en unidades)"
$ws.Cells.Item(274, 15).Value = "Perú"
$ws.Cells.Item(274, 16).Value = 850
$ws.Cells.Item(274, 17).Value = 1
$ws.Cells.Item(274, 18).Value = "Hortaliza"

# Make sure the date cell keeps the date number format like the rest of column D
$ws.Cells.Item(274, 4).NumberFormat = $ws.Cells.Item(275, 4).NumberFormat
